# "start working at Tarife" - add a new "Tarif" worksheet at the end of the
# workbook, populate it with the first rows of the price list, and make it
# the active tab.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing one, so it lands at the end of
# the tab strip (rather than Excel's default of inserting before the active
# sheet).
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($sheetCount))
$ws.Name = "Tarif"

# Seed the sheet's starter content.
$ws.Range("B2").Value = "Abonament sala"
$ws.Range("B3").Value = "Antrenor"
$ws.Range("D2").Value = "dropdown - standard, de familie"
$ws.Range("D3").Value = "dropdown - nume antrenori"
$ws.Range("B4").Value = "Abonament sauna"
$ws.Range("B5").Value = "Sedinta sauna"

# Match the author's last selection on the new (now active) sheet.
$ws.Range("H11").Select() | Out-Null
